$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the Price/Volume columns so that
# numeric-looking strings (e.g. "1.00", "63.627.43") are stored as text,
# matching the source data which uses inline strings, not numbers.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = '63.627.43'
$ws.Range("E2").Value = '  -3.18%  '
$ws.Range("D3").Value = '2.608.38'
$ws.Range("E3").Value = '  -1.99%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '573.55'
$ws.Range("E5").Value = '  -3.93%  '
$ws.Range("D6").Value = '154.91'
$ws.Range("E6").Value = '  -1.28%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -5.12%  '
$ws.Range("E9").Value = '  -6.32%  '
$ws.Range("E10").Value = '  -0.27%  '
$ws.Range("D11").Value = '0.382'
$ws.Range("E11").Value = '  -4.82%  '
$ws.Range("E12").Value = '  -0.53%  '
$ws.Range("D13").Value = '28.24'
$ws.Range("E13").Value = '  -1.81%  '
$ws.Range("D14").Value = '3.077.95'
$ws.Range("E14").Value = '  -2.00%  '
$ws.Range("E15").Value = '  -7.21%  '
$ws.Range("D16").Value = '63.448.26'
$ws.Range("E16").Value = '  -3.24%  '
$ws.Range("D17").Value = '2.605.94'
$ws.Range("E17").Value = '  -2.68%  '
$ws.Range("E18").Value = '  -4.51%  '
$ws.Range("D19").Value = '7.50'
$ws.Range("E19").Value = '  +1.38%  '
$ws.Range("D20").Value = '4.54'
$ws.Range("E20").Value = '  -4.77%  '
$ws.Range("D21").Value = '342.84'
$ws.Range("E21").Value = '  -1.70%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '67.05'
$ws.Range("D24").Value = '1.76'
$ws.Range("E24").Value = '  -2.60%  '
$ws.Range("D25").Value = '0.0000107'
$ws.Range("E25").Value = '  -3.33%  '
$ws.Range("D26").Value = '585.88'
$ws.Range("E26").Value = '  +3.28%  '
$ws.Range("D27").Value = '9.14'
$ws.Range("E27").Value = '  -3.82%  '
$ws.Range("D28").Value = '1.56'
$ws.Range("E28").Value = '  -3.86%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("E30").Value = '  -1.31%  '
$ws.Range("D31").Value = '7.91'
$ws.Range("E31").Value = '  -1.79%  '
$ws.Range("E32").Value = '  -3.23%  '
$ws.Range("E33").Value = '  -4.78%  '
$ws.Range("D34").Value = '6.53'
$ws.Range("E34").Value = '  -2.12%  '
$ws.Range("D35").Value = '5.39'
$ws.Range("E35").Value = '  -0.84%  '
$ws.Range("D36").Value = '0.403'
$ws.Range("E36").Value = '  -4.10%  '
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").Value = '19.64'
$ws.Range("E38").Value = '  -4.14%  '
$ws.Range("D39").Value = '154.08'
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("E40").Value = '  -3.33%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '41.45'
$ws.Range("E42").Value = '  -3.01%  '
$ws.Range("D43").Value = '2.44'
$ws.Range("E43").Value = '  +7.51%  '
$ws.Range("D44").Value = '155.55'
$ws.Range("E44").Value = '  -3.00%  '
$ws.Range("E45").Value = '  -4.65%  '
$ws.Range("D46").Value = '23.22'
$ws.Range("E46").Value = '  +2.28%  '
$ws.Range("D47").Value = '0.0587'
$ws.Range("E47").Value = '  -3.74%  '
$ws.Range("E48").Value = '  -1.90%  '
$ws.Range("D49").Value = '0.1000'
$ws.Range("E49").Value = '  -2.07%  '
$ws.Range("E50").Value = '  -3.34%  '
$ws.Range("D51").Value = '18.78'
$ws.Range("E51").Value = '  -4.91%  '

# Restore the default (Normal) cell style so formatting matches the original
# workbook (only the text content changed, not the cell style).
$priceVolRange.Style = "Normal"
